# edit.ps1
#
# Applies two changes to the document:
#
#   1. The Title paragraph's text runs
#         "Dummy" + " 2" + " " + "Trial " + "2 " + "Test Case Document - ..."
#      are collapsed down to
#         "Edited " + "Test Case Document - ..."
#      (the trailing run is left as its own run, the five leading runs
#      become a single "Edited " run).
#
#   2. The <w:lastRenderedPageBreak/> marker that currently sits on the
#      "Input: abcd@xyz" paragraph of TC004 is moved to the following
#      "Expected Result: Rejected - Must include at least one number"
#      paragraph.
#
# Because plain Range.Text / InsertBefore / Delete edits on this engine
# cause the whole paragraph's runs to be re-merged, we use
# Range.InsertXML with a fully-formed <w:p> fragment (built from the
# paragraph's own original attributes) so the exact run layout called
# for by the diff is produced.

$d = $word.ActiveDocument
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"'

# ---------------------------------------------------------------------
# 1. Title paragraph
# ---------------------------------------------------------------------
$titlePara = $d.Paragraphs(1)

$titleXml = '<w:p ' + $wNs + ' w14:paraId="739F5A3E" w14:textId="010D9CC0" w:rsidR="002E4A0C" w:rsidRDefault="00696F63" w:rsidP="00B56C6E">' +
            '<w:pPr><w:pStyle w:val="Title"/></w:pPr>' +
            '<w:r><w:t xml:space="preserve">Edited </w:t></w:r>' +
            '<w:r><w:t>Test Case Document – Password Policy Including Special Character</w:t></w:r>' +
            '</w:p>'

[void]$titlePara.Range.InsertXML($titleXml)

# ---------------------------------------------------------------------
# 2. Move <w:lastRenderedPageBreak/> from TC004's "Input:" paragraph to
#    its "Expected Result:" paragraph.
# ---------------------------------------------------------------------
$anchor = $d.Content
[void]$anchor.Find.Execute(
    "Description: Password with special character but no number",
    $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$anchor.Collapse(0)

[void]$anchor.Find.Execute(
    "Input: abcd@xyz",
    $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

$inputPara = $anchor.Paragraphs(1)
$expectedPara = $inputPara.Next()

$inputXml = '<w:p ' + $wNs + ' w14:paraId="27E0FDDE" w14:textId="77777777" w:rsidR="002E4A0C" w:rsidRDefault="00000000">' +
            '<w:r><w:t>Input: abcd@xyz</w:t></w:r>' +
            '</w:p>'
[void]$inputPara.Range.InsertXML($inputXml)

$expectedXml = '<w:p ' + $wNs + ' w14:paraId="7D780310" w14:textId="77777777" w:rsidR="002E4A0C" w:rsidRDefault="00000000">' +
               '<w:r><w:lastRenderedPageBreak/><w:t>Expected Result: Rejected – Must include at least one number</w:t></w:r>' +
               '</w:p>'
[void]$expectedPara.Range.InsertXML($expectedXml)

Write-Host "Edits applied."
